$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price (D) and volume-change (E) values.
# D-column values are forced to remain plain text (matching the original
# inline-string storage) even when they look numeric, by temporarily
# applying a text number format and then restoring the default style so
# no extra cell formatting is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.601.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.725.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.00%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.53%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.547"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.725.09"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.98%  "
$ws.Range("E10").Value = "  +2.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.366"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.37%  "
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.225.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.98%  "
$ws.Range("E16").Value = "  +1.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.673.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.720.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "374.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.36%  "
$ws.Range("E22").Value = "  +2.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.75%  "
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.867.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000105"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "589.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("E32").Value = "  +3.29%  "
$ws.Range("E33").Value = "  +3.58%  "
$ws.Range("E34").Value = "  +5.83%  "
$ws.Range("E35").Value = "  +2.80%  "
$ws.Range("E36").Value = "  -1.72%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "162.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.46%  "
$ws.Range("E40").Value = "  +2.80%  "
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("E42").Value = "  +2.62%  "
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0311"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.33%  "
$ws.Range("E48").Value = "  +4.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "155.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("E50").Value = "  +3.35%  "
$ws.Range("E51").Value = "  +5.64%  "
